$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 94, shifting existing rows (94-107) down to (95-108).
$ws.Rows.Item(94).Insert()

# Copy the date number format of column D from the row pushed down (now row 95) to the new row 94.
$ws.Cells.Item(94, 4).NumberFormat = $ws.Cells.Item(95, 4).NumberFormat

# Fill constant columns (same across all records in this block) for the new row 94.
$ws.Cells.Item(94, 1).Value = 10
$ws.Cells.Item(94, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(94, 3).Value = "La Araucanía"
$ws.Cells.Item(94, 4).Value = 45127
$ws.Cells.Item(94, 5).Value = 9
$ws.Cells.Item(94, 6).Value = "Fruta"
$ws.Cells.Item(94, 7).Value = 100108
$ws.Cells.Item(94, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(94, 9).Value = 100108003
$ws.Cells.Item(94, 10).Value = "Maracuyá"
$ws.Cells.Item(94, 11).Value = "Sin especificar"
$ws.Cells.Item(94, 12).Value = "Primera"
$ws.Cells.Item(94, 13).Value = 15
$ws.Cells.Item(94, 14).Value = 45000
$ws.Cells.Item(94, 15).Value = 45000
$ws.Cells.Item(94, 16).Value = 45000
$ws.Cells.Item(94, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(94, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(94, 19).Value = 2500
$ws.Cells.Item(94, 20).Value = 18
